$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# 1. Sheet "Metadata": update Version + Date values
# -----------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value2 = "2.0.1-sd-202510-matchbox-patch"
$meta.Range("B8").Value2 = "2025-10-29T22:15:57+01:00"

# -----------------------------------------------------------------------
# 2. Sheet "Metadata": insert a new "Jurisdiction" property row right
#    after "Contact" (row 10) and before "Description" (row 11).
#    Shift rows 11..19 down to 12..20 first (bottom-up so we never
#    clobber a row before it has been read), then write the new row.
# -----------------------------------------------------------------------
for ($r = 19; $r -ge 11; $r--) {
    $n = $r + 1
    $meta.Range("A" + $n + ":B" + $n).ClearContents()
    $meta.Range("A" + $r + ":B" + $r).Copy($meta.Range("A" + $n + ":B" + $n))
}

$meta.Range("A11").Value2 = "Jurisdiction"
$meta.Range("B11").Value2 = ""

# -----------------------------------------------------------------------
# 3. Sheet "Elements": add the II-1 invariant constraint text to the
#    "RecordTarget.typeId" row (row 5), column AJ ("Constraint(s)").
# -----------------------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AJ5").Value2 = "II-1:An II instance must have either a root or an nullFlavor. {root.exists() or nullFlavor.exists()}`n"
